$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (46081 -> 46082) for every data row (rows 2 through 273).
$ws.Range("C2:C273").Value = 46082
